$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.256.32'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.830.00'
$ws.Range('E3').Value = '  +2.65%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.560'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.07'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +3.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0724'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0931'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '2.093.02'
$ws.Range('E12').Value = '  +2.74%  '
$ws.Range('D13').Value = '1.829.97'
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.646'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('D16').Value = '34.276.49'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').Value = '0.0₃0793'
$ws.Range('E20').Value = '  +7.14%  '
$ws.Range('E21').Value = '  +8.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.30'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.94%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.71%  '
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0535'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.21%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.90'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.81%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.441.03'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.649'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.91%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.60%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0190'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.970'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.85%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '82.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.20%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.84%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.62%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.990.02'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0497'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.91%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.05'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '107.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.05%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.92'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.96%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0125'
$ws.Range('E51').Value = '  +5.84%  '
